## Generate Report for Archive
##
## The localization status report is regenerated: the "Status" value for the
## tracked file moves from "Ready for handoff" to "In Translation" (Overview
## sheet columns E/F which mirror the per-locale Status column, plus the
## Status column on each locale sheet). Regenerating the report also re-runs
## the column autosizing for those Status columns, so their width shrinks to
## fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears ------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2. Re-fit the Status columns that shrank as a result --------------
# ColumnWidth is expressed in characters; 12.5 characters is the closest
# representable width to the regenerated report's target column width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
